$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.266.65"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "2.268.57"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'305.41"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'97.04"
$ws.Range("E6").Value = "  +4.50%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").Value = "'35.54"
$ws.Range("E10").Value = "  +8.69%  "

$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "'6.65"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "2.616.98"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "'14.36"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "2.272.71"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").Value = "'0.794"
$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").Value = "42.157.87"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").Value = "'12.48"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "'5.96"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'67.53"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").Value = "'237.28"
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'1.96"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "'23.80"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "'37.44"
$ws.Range("E28").Value = "  +6.05%  "

$ws.Range("D29").Value = "'9.50"
$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "'160.04"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").Value = "'5.25"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "'3.15"
$ws.Range("E34").Value = "  +4.49%  "

$ws.Range("D35").Value = "'0.0741"
$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("D36").Value = "'17.09"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").Value = "'1.83"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("D41").Value = "'4.06"
$ws.Range("E41").Value = "  +3.29%  "

$ws.Range("D42").Value = "'2.43"
$ws.Range("E42").Value = "  +14.27%  "

$ws.Range("D43").Value = "1.992.15"
$ws.Range("E43").Value = "  -0.59%  "

$ws.Range("D44").Value = "'0.0286"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").Value = "'18.75"
$ws.Range("E45").Value = "  -4.67%  "

$ws.Range("D46").Value = "'9.97"
$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("D48").Value = "'53.23"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'72.05"
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").Value = "'91.34"
$ws.Range("E51").Value = "  +0.03%  "
